# SRG_YR_FIN.xlsx update — adds the latest reporting-period column to the
# three financial statement tables (Income Statement, Balance Sheet, Cash
# Flow Statement) on sheet "SRG" by inserting a new column D and filling in
# the newest period's figures, shifting the previously existing periods one
# column to the right (D->E, E->F, ... K->L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D - this shifts the existing
# D:K data (and any per-cell styles) right to E:L automatically.
$ws.Range("D1").EntireColumn.Insert()

# The freshly inserted column D cells don't carry the row's established
# number formatting yet (Insert() pulls formatting from the column to the
# left instead). Re-apply the correct formatting by copying it over from
# column E (which now holds what used to be column D), for every data row.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Income Statement (rows 7-35) -----------------------------------------
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 214800
$ws.Range("D9").Value = 28700
$ws.Range("D10").Value = 186000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 226700
$ws.Range("D17").Value = 332900
$ws.Range("D18").Value = -118100
$ws.Range("D20").Value = 93600
$ws.Range("D21").Value = 202100
$ws.Range("D22").Value = 90000
$ws.Range("D23").Value = -114600
$ws.Range("D24").Value = 300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -114900
$ws.Range("D27").Value = -78400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -93600
$ws.Range("D33").Value = -78400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -78400

# --- Balance Sheet (rows 38-77) --------------------------------------------
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 532900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 36900
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 398600
$ws.Range("D48").Value = 1751100
$ws.Range("D49").Value = 123700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 3100
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2876100
$ws.Range("D57").Value = 71000
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 34800
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 1598100
$ws.Range("D62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2095300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -344100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 780800
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (rows 80-102) -------------------------------------
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -78400
$ws.Range("D83").Value = 226700
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 54900
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -119500
$ws.Range("D96").Value = -39700
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 180200
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 115600

# Row 12 is one of the few rows that didn't have numeric historical data
# (was already "NA" across D:K) - the shift leaves it "NA" too.
$ws.Range("D12").Value = "NA"
